$d = $word.ActiveDocument

# Several merge-field placeholders in this template were split across multiple
# <w:r> runs, e.g. "${" + "code" + "}" instead of a single "${code}" run.
# Collapse each split placeholder back into a single run while preserving the
# character formatting of the run that carries the field name itself.
function Fix-SplitPlaceholder($fieldName) {
    $pattern = '${' + $fieldName + '}'

    # Locate the placeholder text (it may currently be spread across runs;
    # Word's Find still matches the concatenated story text).
    $locate = $d.Content.Duplicate
    $found = $locate.Find.Execute($pattern, $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if (-not $found) {
        return
    }

    # Strip the leading "${" run-text, scoped to just this occurrence so the
    # other placeholders elsewhere in the document are left untouched.
    $target1 = $d.Range($locate.Start, $locate.End)
    $target1.Find.Execute('${', $false, $false, $false, $false, $false, $true, 0, $false, '', 1) | Out-Null

    # Strip the trailing "}" run-text, again scoped to this occurrence only.
    $target2 = $d.Range($locate.Start, $locate.End)
    $target2.Find.Execute('}', $false, $false, $false, $false, $false, $true, 0, $false, '', 1) | Out-Null

    # What remains is a single run containing just the field name; expand it
    # back out to the full "${fieldName}" placeholder text.
    $target3 = $d.Range($locate.Start, $locate.End)
    $target3.Find.Execute($fieldName, $false, $false, $false, $false, $false, $true, 0, $false, $pattern, 1) | Out-Null
}

Fix-SplitPlaceholder('code')
Fix-SplitPlaceholder('national')
Fix-SplitPlaceholder('school_date')
Fix-SplitPlaceholder('category')

# The job level field was bound to the wrong property; point it at the name
# instead of the id.
$d.Content.Find.Execute('${job_level_id}', $false, $false, $false, $false, $false, $true, 1, $false, '${job_level_name}', 2) | Out-Null
